$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 02:16"

# --- China (row 4) updated case numbers ---
$ws.Range("B4").Value = 81171
$ws.Range("C4").Value = 78
$ws.Range("D4").Value = 73159
$ws.Range("E4").Value = 4735
$ws.Range("F4").Value = 1573
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 3277

# --- Estados Unidos (row 6) updated case numbers ---
$ws.Range("B6").Value = 43721
$ws.Range("C6").Value = 10155
$ws.Range("G6").Value = 139
$ws.Range("H6").Value = 552

# --- Venezuela moves up (now row 94), with refreshed data; Bielorrusia and
#     Senegal shift down one row each (rows 95 and 96), keeping their
#     previous case numbers ---
$ws.Range("A94").Value = "Venezuela"
$ws.Range("B94").Value = 84
$ws.Range("C94").Value = 14
$ws.Range("D94").Value = 15
$ws.Range("E94").Value = 69
$ws.Range("F94").Value = 2
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 0

$ws.Range("A95").Value = "Bielorrusia"
$ws.Range("B95").Value = 81
$ws.Range("C95").Value = 5
$ws.Range("D95").Value = 22
$ws.Range("E95").Value = 59
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0

$ws.Range("A96").Value = "Senegal"
$ws.Range("B96").Value = 79
$ws.Range("C96").Value = 12
$ws.Range("D96").Value = 8
$ws.Range("E96").Value = 71
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
